# Refresh the crypto price/volume columns (D = Price, E = Volume(1h))
# with the latest scrape. Values are stored as text (matching the sheet's
# existing inline-string cells), so each cell is forced to Text format
# before the write to stop Excel from auto-converting the numeric- and
# percent-looking strings into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.96%"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.68%"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.152"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.07%"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07384"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.20%"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.827"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "24.50%"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.854"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.53%"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.751"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.83%"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9292"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.78%"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1704"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.37%"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07173"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.67%"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08060"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.36%"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03030"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.86%"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09941"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.40%"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001498"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.83%"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006102"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-5.71%"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.469"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.30%"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.224"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.23%"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3257"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.37%"
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.21%"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.585"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.53%"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04644"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.57%"
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.71%"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001214"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.46%"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004768"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.87%"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001296"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.49%"
# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "7.27%"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01731"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.47%"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04519"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.26%"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007113"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.39%"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1344"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.14%"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002164"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.17%"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01098"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-13.41%"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006222"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.11%"
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-21.60%"
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7399"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "4.33%"
